# Apply the new built-in Table Style to the table on slide 6.
#
# The source XML shows the table's <a:tblPr><a:tableStyleId> changing from
#   {F90B60FD-4CB7-4C4D-9432-BE31291BA696}
# to
#   {7D373626-7F27-4839-B458-CAA57E3A5D15}
# which is exactly what happens when the user selects the table and picks a
# different style from the PowerPoint "Table Styles" gallery (Table Design
# tab). In the PowerPoint object model this is done with Table.ApplyStyle.

$p = $ppt.ActivePresentation

$newStyleId = "{7D373626-7F27-4839-B458-CAA57E3A5D15}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
